$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.228.15"
$ws.Range("E2").Value = "  -0.67%  "
$ws.Range("D3").Value = "1.858.40"
$ws.Range("E3").Value = "  -1.18%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.000"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -0.11%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "242.10"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.63%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.6954"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -2.98%  "
$ws.Range("E7").Value = "  -0.03%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.07808"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -1.89%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.3118"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -0.97%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "24.04"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -3.83%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.07813"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -4.20%  "
$ws.Range("D12").Value = "1.877.31"
$ws.Range("E12").Value = "  -0.59%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "5.135"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -2.15%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "92.23"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -2.69%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.6932"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -2.32%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "6.548"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +2.21%  "
$ws.Range("E17").Value = "  +0.28%  "
$ws.Range("D18").Value = "29.272.18"
$ws.Range("E18").Value = "  -0.59%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "248.83"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -2.13%  "
$ws.Range("D20").Value = "2.107.42"
$ws.Range("E20").Value = "  -2.11%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "12.92"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -3.12%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.9997"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.12%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "7.554"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -2.41%  "
$ws.Range("E24").Value = "  -0.30%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "0.1529"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -3.49%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "161.03"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -0.93%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "8.912"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -1.89%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "18.62"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -1.88%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "1.576"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +4.62%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "4.273"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -3.34%  "
$ws.Range("E31").Value = "  -1.13%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "1.203"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -2.04%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.05228"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -1.88%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "1.883"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -3.43%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.7570"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -0.18%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "1.173"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -0.54%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "2.696"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -0.07%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.01859"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -2.14%  "
$ws.Range("D39").Value = "1.241.37"
$ws.Range("E39").Value = "  -2.77%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "2.738"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -0.99%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.9014"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -0.74%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "110.97"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -1.83%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "5.909"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -7.91%  "
$ws.Range("E44").Value = "  +0.00%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "68.65"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -7.68%  "
$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").Value = "2.006.36"
$ws.Range("E46").Value = "  -1.65%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "9.574"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +0.87%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.5178"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -0.55%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.00000000122"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -6.72%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "1.773"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -2.05%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.4261"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -2.17%  "
